$p = $ppt.ActivePresentation

# Remove slide 2 (the "Trend plot / Sensitivity plot" slide) from the deck.
$p.Slides.Item(2).Delete()
